# The deck ships two theme parts: ppt/theme/theme1.xml ("Office Theme") and
# ppt/theme/theme2.xml ("Integral") -- the slide master / notes master /
# presentation currently resolve to the "Integral" palette. The edit swaps
# the two themes so the deck uses the plain "Office Theme" palette instead.
#
# The PowerPoint object model's theme-color surface (Master.ColorScheme /
# NotesMaster.ColorScheme / Slide.ColorScheme, etc.) all resolve to the same
# single live theme part, so we drive the swap through that one shared
# ColorScheme object, writing each of the 12 scheme slots (dk1, lt1, dk2,
# lt2, accent1-6, hlink, folHlink) with the "Office Theme" RGB values in
# the same slot order used by OOXML's <a:clrScheme>.

function HexToRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# "Office Theme" clrScheme values (currently ppt/theme/theme1.xml), in
# <a:clrScheme> child order.
$officeThemeColors = @(
    "000000",  # dk1
    "FFFFFF",  # lt1
    "44546A",  # dk2
    "E7E6E6",  # lt2
    "5B9BD5",  # accent1
    "ED7D31",  # accent2
    "A5A5A5",  # accent3
    "FFC000",  # accent4
    "4472C4",  # accent5
    "70AD47",  # accent6
    "0563C1",  # hlink
    "954F72"   # folHlink
)

$p = $ppt.ActivePresentation
$colorScheme = $p.SlideMaster.ColorScheme

for ($i = 0; $i -lt $officeThemeColors.Count; $i++) {
    $colorScheme.Colors($i + 1).RGB = HexToRgb $officeThemeColors[$i]
}
